$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns
$ws.Range("T1").Value = "RD_X"
$ws.Range("U1").Value = "RD_Y"

# RD_X / RD_Y values per row (row number => x, y)
$rdValues = @{
    2  = @(121616, 488227)
    3  = @(121613, 488215)
    4  = @(121558, 488127)
    5  = @(121558, 488117)
    6  = @(121599, 488078)
    7  = @(121608, 488071)
    8  = @(121643, 488043)
    9  = @(121657, 488033)
    10 = @(121667, 488026)
    11 = @(121753, 487951)
    12 = @(121761, 487939)
    13 = @(121819, 487739)
    14 = @(121812, 487742)
    15 = @(121806, 487749)
    16 = @(121815, 487751)
    17 = @(121825, 487752)
    18 = @(121832, 487753)
    19 = @(121841, 487754)
    20 = @(121853, 487755)
}

foreach ($row in $rdValues.Keys) {
    $vals = $rdValues[$row]
    $ws.Cells.Item($row, 20).Value = $vals[0]
    $ws.Cells.Item($row, 21).Value = $vals[1]
}

# Update selection to reflect the new active cell shown in the diff
$ws.Range("U20").Select()
